# Update countries & provincias Spain
# Applies refreshed COVID data values and reorders a few country rows
# to match the newly re-sorted dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# List of row-level changes: country name (column A) plus the new
# Casos totales / Nuevos casos / Casos activos / Recuperados /
# Casos criticos / Muertes hoy / Muertes values (columns B..H).
$changes = @(
    @{Row=4;   A="Estados Unidos";                   B=3479650; C=167; D=1549624; E=1791779; F=0; G=0;   H=138247}
    @{Row=6;   A="India";                             B=911606;  C=3961; D=573283; E=314544;  F=0; G=52;  H=23779}
    @{Row=14;  A="Iran";                              B=262173;  C=2521; D=225270; E=23692;   F=0; G=179; H=13211}
    @{Row=36;  A="Filipinas";                         B=57545;   C=634;  D=20459;  E=35483;   F=0; G=6;   H=1603}
    @{Row=37;  A="Kuwait";                             B=56174;  C=666;  D=46161;  E=9617;    F=0; G=3;   H=396}
    @{Row=48;  A="Afganistan";                        B=34740;   C=285;  D=21454;  E=12238;   F=0; G=36;  H=1048}
    @{Row=49;  A="Rumania";                           B=33585;   C=637;  D=21803;  E=9851;    F=0; G=30;  H=1931}
    @{Row=50;  A="Barein";                            B=33476;   C=0;    D=29099;  E=4268;    F=0; G=0;   H=109}
    @{Row=51;  A="Nigeria";                           B=33153;   C=0;    D=13671;  E=18738;   F=0; G=0;   H=744}
    @{Row=65;  A="Marruecos";                         B=16047;   C=111;  D=13403;  E=2388;    F=0; G=1;   H=256}
    @{Row=79;  A="Malasia";                           B=8729;    C=4;    D=8524;   E=83;      F=0; G=0;   H=122}
    @{Row=82;  A="Consejo Danes para los Refugiados"; B=8135;    C=60;   D=3948;   E=3997;    F=0; G=0;   H=190}
    @{Row=86;  A="Finlandia";                         B=7301;    C=6;    D=6800;   E=172;     F=0; G=0;   H=329}
    @{Row=88;  A="Estado de Palestina";               B=6764;    C=198;  D=1084;   E=5638;    F=0; G=3;   H=42}
    @{Row=89;  A="Haiti";                             B=6727;    C=0;    D=3022;   E=3564;    F=0; G=2;   H=141}
    @{Row=90;  A="Tayikistan";                        B=6596;    C=0;    D=5278;   E=1263;    F=0; G=0;   H=55}
    @{Row=102; A="Albania";                           B=3667;    C=96;   D=2062;   E=1508;    F=0; G=2;   H=97}
    @{Row=110; A="Sri Lanka";                         B=2649;    C=3;    D=1988;   E=650;     F=0; G=0;   H=11}
    @{Row=126; A="Hong Kong";                         B=1570;    C=48;   D=1229;   E=333;     F=0; G=0;   H=8}
    @{Row=127; A="Nueva Zelanda";                     B=1545;    C=1;    D=1498;   E=25;      F=0; G=0;   H=22}
    @{Row=139; A="Uganda";                            B=1040;    C=11;   D=984;    E=56;      F=0; G=0;   H=0}
    @{Row=140; A="Burkina Faso";                      B=1036;    C=0;    D=869;    E=114;     F=0; G=0;   H=53}
    @{Row=141; A="Zimbabue";                          B=1034;    C=0;    D=343;    E=672;     F=0; G=0;   H=19}
)

foreach ($change in $changes) {
    $r = $change.Row
    $ws.Cells.Item($r, 1).Value = $change.A
    $ws.Cells.Item($r, 2).Value = $change.B
    $ws.Cells.Item($r, 3).Value = $change.C
    $ws.Cells.Item($r, 4).Value = $change.D
    $ws.Cells.Item($r, 5).Value = $change.E
    $ws.Cells.Item($r, 6).Value = $change.F
    $ws.Cells.Item($r, 7).Value = $change.G
    $ws.Cells.Item($r, 8).Value = $change.H
}

# Update the "last updated" timestamp footer line.
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 12:29"
